# Refined metadata to be additional tab

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# Update time_taken column (F) timestamps on the data sheet
$newTimes = @(
    "2021-10-05 14:35:14.665972",
    "2021-10-05 14:35:14.665980",
    "2021-10-05 14:35:14.665983",
    "2021-10-05 14:35:14.665986",
    "2021-10-05 14:35:14.665989",
    "2021-10-05 14:35:14.665992",
    "2021-10-05 14:35:14.665994",
    "2021-10-05 14:35:14.665997",
    "2021-10-05 14:35:14.666000",
    "2021-10-05 14:35:14.666002",
    "2021-10-05 14:35:14.666005",
    "2021-10-05 14:35:14.666007",
    "2021-10-05 14:35:14.666010",
    "2021-10-05 14:35:14.666013",
    "2021-10-05 14:35:14.666015",
    "2021-10-05 14:35:14.666018",
    "2021-10-05 14:35:14.666020",
    "2021-10-05 14:35:14.666023",
    "2021-10-05 14:35:14.666026",
    "2021-10-05 14:35:14.666028",
    "2021-10-05 14:35:14.666031",
    "2021-10-05 14:35:14.666033",
    "2021-10-05 14:35:14.666036",
    "2021-10-05 14:35:14.666039",
    "2021-10-05 14:35:14.666042"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# Add a new "metadata" worksheet after the "data" sheet
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Match the page margins used by the "data" sheet (points: 1in=72pt)
$metaSheet.PageSetup.LeftMargin = 54
$metaSheet.PageSetup.RightMargin = 54
$metaSheet.PageSetup.TopMargin = 72
$metaSheet.PageSetup.BottomMargin = 72
$metaSheet.PageSetup.HeaderMargin = 36
$metaSheet.PageSetup.FooterMargin = 36

# Header row
$metaSheet.Cells.Item(1, 2).Value = "data_name"
$metaSheet.Cells.Item(1, 3).Value = "data_id"
$metaSheet.Cells.Item(1, 4).Value = "data_version"
$metaSheet.Cells.Item(1, 5).Value = "data_version_created"
$metaSheet.Cells.Item(1, 6).Value = "panel_query_time"
$metaSheet.Cells.Item(1, 7).Value = "panel_get_request"

# Data row
$metaSheet.Cells.Item(2, 1).Value = 0
$metaSheet.Cells.Item(2, 2).Value = "Photosensitivity Syndromes"
$metaSheet.Cells.Item(2, 3).Value = 156
$metaSheet.Cells.Item(2, 4).NumberFormat = "@"
$metaSheet.Cells.Item(2, 4).Value = "1.0"
$metaSheet.Cells.Item(2, 4).Style = "Normal"
$metaSheet.Cells.Item(2, 5).Value = "2021-04-24T01:18:47.889043Z"
$metaSheet.Cells.Item(2, 6).Value = "2021-10-05 14:35:14.662259"
$metaSheet.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/156/?format=json"

# Copy the header style (bold + border, matching the "data" sheet header) onto
# the new header row (B1:G1) and onto A2 (row index column), mirroring the
# formatting already used on the "data" sheet.
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)

$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

$metaSheet.Select()
$metaSheet.Range("A1").Select()
